$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Add the new "Sheet2" (holds the list of "Loại điều chuyển" options) right after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"
$ws2.Range("A1").Value = "Chuyển vị trí"
$ws2.Range("A2").Value = "Nghỉ việc"
$ws2.Range("A3").Value = "Nghỉ thai sản"
$ws2.Range("A4").Value = "Thai sản đi làm lại"

# New data row on Sheet1
$ws1.Range("A3").Value = 12543
$ws1.Range("D3").Value = "Nghỉ việc"
$ws1.Range("E3").Value = [DateTime]"2024-08-15"
$ws1.Range("F3").Value = "Rút hồ sơ"

# Turn the "Loại điều chuyển" column into a picklist instead of free text entry
$rng = $ws1.Range("D2:D70")
$rng.Validation.Add(3, 1, 1, "=Sheet2!`$A`$1:`$A`$4")

# Keep the active selection in sync with the new last row
$ws1.Range("F3").Select() | Out-Null
